$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '72.328.98'
$ws.Range("E2").Value = '  +3.90%  '
# Row 3
$ws.Range("D3").Value = '4.021.75'
$ws.Range("E3").Value = '  +2.73%  '
# Row 4
$ws.Range("E4").Value = '  +0.14%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '514.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.21%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.52%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.714'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +16.33%  '
# Row 8
$ws.Range("E8").Value = '  +0.12%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.764'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.75%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.173'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.18%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000322'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.60%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.49'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +10.13%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.75'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.59%  '
# Row 14
$ws.Range("D14").Value = '4.664.50'
$ws.Range("E14").Value = '  +3.12%  '
# Row 15
$ws.Range("D15").Value = '4.026.31'
$ws.Range("E15").Value = '  +2.41%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.55%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.64%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.52%  '
# Row 19
$ws.Range("E19").Value = '  -1.78%  '
# Row 20
$ws.Range("D20").Value = '72.157.48'
$ws.Range("E20").Value = '  +4.02%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '434.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.48%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '102.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +15.92%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.09%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.59'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.40%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.47%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.56'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.66%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.84%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.92%  '
# Row 29
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.79'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.09%  '
# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.04%  '
# Row 31
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.51%  '
# Row 32
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '677.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.57%  '
# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.127'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.94%  '
# Row 34
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.79'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +14.36%  '
# Row 35
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '67.87'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.96%  '
# Row 36
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '40.82'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.68%  '
# Row 37
$ws.Range("D37").Value = '0.0₃0865'
$ws.Range("E37").Value = '  +2.59%  '
# Row 38
$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.431'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.82%  '
# Row 39
$ws.Range("B39").Value = 'ThetaToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.50'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +14.80%  '
# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.151'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.42%  '
# Row 41
$ws.Range("B41").Value = 'Dai'
$ws.Range("C41").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.05%  '
# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.11%  '
# Row 43
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0485'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.91%  '
# Row 44
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.41%  '
# Row 45
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.157'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.48%  '
# Row 46
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.90%  '
# Row 47
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.39%  '
# Row 48
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.78%  '
# Row 49
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.12%  '
# Row 50
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000269'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +18.06%  '
# Row 51
$ws.Range("B51").Value = 'LidoDAOToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.48%  '
